# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# --- 1. Insert a new "2022-Q1" sheet right before "总计" -------------------
# NOTE: Worksheets.Add() always inserts the new sheet at the front of the
# tab strip, and Move() repositions by index too - both operations shift
# every other sheet's index. Any sheet reference obtained before such a
# call goes stale, so we always (re)fetch sheets by name right before we
# use them.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"
$newSheet.Move($wb.Worksheets.Item("总计"))
$newSheet = $wb.Worksheets.Item("2022-Q1")

# Reference sheet used purely to clone the existing "bold / centered /
# bordered" header-row & index-column look (style index shared by the
# other quarterly sheets) onto the freshly-added sheet.
$styleSource = $wb.Worksheets.Item("2021-Q4")

# --- 2. Populate "2022-Q1" with the fund holdings table --------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @(0, "519995", "长信金利趋势混合",                                 "56.80", "86.00", "2.96", "1.6813", 6),
    @(1, "159851", "华宝中证金融科技主题ETF",                          "3.16",  "98.58", "4.26", "0.1346", 5),
    @(2, "010153", "中加中证500指数增强A",                             "1.26",  "94.19", "3.33", "0.0420", 3),
    @(3, "516100", "华夏中证金融科技主题交易型开放式指数证券投资基金", "0.68",  "96.91", "4.22", "0.0287", 5),
    @(4, "010154", "中加中证500指数增强C",                             "0.60",  "94.19", "3.33", "0.0200", 3)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").NumberFormat = "@"
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Clone the bold/centered/bordered look of the header row + index column
# from an existing sheet so the new tab visually matches its siblings.
$styleSource.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)
$styleSource.Range("A2:A6").Copy()
$newSheet.Range("A2:A6").PasteSpecial($xlPasteFormats)

# --- 3. Insert a new top data row into "总计" for 2022-Q1 -------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows("2").Insert()
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 1.91

# Clone A3's index-column style onto the newly-inserted A2 so it matches
# the rest of the column ("总计" already uses that look for A3:A5).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)
